# Update "想去人数" (want-to-go count) figures in column F across the
# four sheets of 北京-漫展信息.xlsx to match a newer data pull.
#
# Sheet order (per workbook.xml): 1=展览, 2=演出, 3=本地生活, 4=全部类型
# (全部类型 is a rollup of the other three sheets, so the same events'
# counts are updated there too.)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$sheet1Updates = @{
    3  = 536
    4  = 686
    7  = 1159
    10 = 2077
    13 = 33
    14 = 63
    15 = 516
    20 = 399
    21 = 399
    22 = 746
    23 = 453
    24 = 2884
    26 = 113
    27 = 3212
    28 = 666
    29 = 535
    30 = 242
    31 = 987
    34 = 721
    35 = 703
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$sheet2Updates = @{
    10 = 65
    13 = 81
    21 = 191
    22 = 136
    23 = 449
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$sheet3Updates = @{
    3 = 2936
    4 = 378
    5 = 254
    6 = 398
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3Updates[$row]
}

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$sheet4Updates = @{
    4  = 536
    6  = 378
    7  = 254
    8  = 686
    12 = 1159
    15 = 398
    16 = 2077
    19 = 33
    20 = 65
    22 = 63
    23 = 516
    25 = 81
    31 = 399
    32 = 399
    34 = 746
    35 = 453
    37 = 2884
    38 = 113
    39 = 3212
    40 = 666
    41 = 535
    42 = 242
    43 = 987
    46 = 191
    47 = 136
    48 = 449
    50 = 721
    51 = 703
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
